# Increase MaxInvest Storage Adapt Szenarios Existing Units
# Column S on the "Power Storage" sheet is "MaxInvest" ([0-N]); raise it
# from 8 to 15 MW for the five existing storage units (rows 7-11).

$wb = $excel.ActiveWorkbook
$ws = $wb.Sheets.Item("Power Storage")

$ws.Range("S7:S11").Value = 15

# Leave the selection on the just-edited range, matching the saved view.
[void]$ws.Range("S8:S11").Select()
